$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 15: remove the "Spreadsheet SpreadsheetResult calc()" header cell entirely (row becomes empty)
$ws.Range("B15").Clear()

# Row 16: clear the "Step Name" / "Value" header cells
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()

# Row 17: clear Step1 and its formula/result value (keep the styled C17 cell as an empty placeholder)
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

# Rows 18-22: clear the Step labels (B column) entirely (format + content) and the
# corresponding result values in column C (keep their style, clear their content)
$ws.Range("B18").Clear()
$ws.Range("C18").ClearContents()

$ws.Range("B19").Clear()
$ws.Range("C19").ClearContents()

$ws.Range("B20").Clear()
$ws.Range("C20").ClearContents()

$ws.Range("B21").Clear()
$ws.Range("C21").ClearContents()

$ws.Range("B22").Clear()
$ws.Range("C22").ClearContents()

# Update the selection to match the reverted view state
$ws.Range("A15:I29").Select() | Out-Null

Write-Host "Edit complete"
